$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Wnt9a"
$ws.Cells.Item(2, 3).Value = "Fzd4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.232908333333333
$ws.Cells.Item(2, 8).Value = 9.698725
$ws.Cells.Item(2, 9).Value = 0.5584514397475191
$ws.Cells.Item(2, 10).Value = 0.5584514397475192
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 19.48350866666667
$ws.Cells.Item(2, 14).Value = 58.450526
$ws.Cells.Item(2, 15).Value = 0.3081250754721727
$ws.Cells.Item(2, 16).Value = 0.3081250754721726
$ws.Cells.Item(2, 17).Value = 62.98839753103888
$ws.Cells.Item(2, 18).Value = 566.89557777935
$ws.Cells.Item(2, 19).Value = 0.1720728920197478
$ws.Cells.Item(2, 20).Value = 0.1720728920197478

$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Wnt9a"
$ws.Cells.Item(3, 3).Value = "Fzd4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.232908333333333
$ws.Cells.Item(3, 8).Value = 9.698725
$ws.Cells.Item(3, 9).Value = 0.5584514397475191
$ws.Cells.Item(3, 10).Value = 0.5584514397475192
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 20.343383
$ws.Cells.Item(3, 14).Value = 61.03014900000001
$ws.Cells.Item(3, 15).Value = 0.3217236961512193
$ws.Cells.Item(3, 16).Value = 0.3217236961512193
$ws.Cells.Item(3, 17).Value = 65.76829242889167
$ws.Cells.Item(3, 18).Value = 591.9146318600251
$ws.Cells.Item(3, 19).Value = 0.1796670613165418
$ws.Cells.Item(3, 20).Value = 0.1796670613165418

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Wnt9a"
$ws.Cells.Item(4, 3).Value = "Fzd4"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.232908333333333
$ws.Cells.Item(4, 8).Value = 9.698725
$ws.Cells.Item(4, 9).Value = 0.5584514397475191
$ws.Cells.Item(4, 10).Value = 0.5584514397475192
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.1305583333333333
$ws.Cells.Item(4, 14).Value = 0.391675
$ws.Cells.Item(4, 15).Value = 0.002064735720865253
$ws.Cells.Item(4, 16).Value = 0.002064735720865253
$ws.Cells.Item(4, 17).Value = 0.4220831238194444
$ws.Cells.Item(4, 18).Value = 3.798748114375
$ws.Cells.Item(4, 19).Value = 0.001153054636015332
$ws.Cells.Item(4, 20).Value = 0.001153054636015333

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Wnt9a"
$ws.Cells.Item(5, 3).Value = "Fzd4"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.232908333333333
$ws.Cells.Item(5, 8).Value = 9.698725
$ws.Cells.Item(5, 9).Value = 0.5584514397475191
$ws.Cells.Item(5, 10).Value = 0.5584514397475192
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 23.275017
$ws.Cells.Item(5, 14).Value = 69.825051
$ws.Cells.Item(5, 15).Value = 0.3680864926557428
$ws.Cells.Item(5, 16).Value = 0.3680864926557428
$ws.Cells.Item(5, 17).Value = 75.245996417775
$ws.Cells.Item(5, 18).Value = 677.2139677599749
$ws.Cells.Item(5, 19).Value = 0.2055584317752142
$ws.Cells.Item(5, 20).Value = 0.2055584317752142

$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Wnt9a"
$ws.Cells.Item(6, 3).Value = "Fzd4"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.556150666666667
$ws.Cells.Item(6, 8).Value = 7.668452
$ws.Cells.Item(6, 9).Value = 0.4415485602524809
$ws.Cells.Item(6, 10).Value = 0.441548560252481
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 19.48350866666667
$ws.Cells.Item(6, 14).Value = 58.450526
$ws.Cells.Item(6, 15).Value = 0.3081250754721727
$ws.Cells.Item(6, 16).Value = 0.3081250754721726
$ws.Cells.Item(6, 17).Value = 49.80278366730578
$ws.Cells.Item(6, 18).Value = 448.225053005752
$ws.Cells.Item(6, 19).Value = 0.1360521834524249
$ws.Cells.Item(6, 20).Value = 0.1360521834524248

$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Wnt9a"
$ws.Cells.Item(7, 3).Value = "Fzd4"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.556150666666667
$ws.Cells.Item(7, 8).Value = 7.668452
$ws.Cells.Item(7, 9).Value = 0.4415485602524809
$ws.Cells.Item(7, 10).Value = 0.441548560252481
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 20.343383
$ws.Cells.Item(7, 14).Value = 61.03014900000001
$ws.Cells.Item(7, 15).Value = 0.3217236961512193
$ws.Cells.Item(7, 16).Value = 0.3217236961512193
$ws.Cells.Item(7, 17).Value = 52.00075201770535
$ws.Cells.Item(7, 18).Value = 468.0067681593481
$ws.Cells.Item(7, 19).Value = 0.1420566348346775
$ws.Cells.Item(7, 20).Value = 0.1420566348346775

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Wnt9a"
$ws.Cells.Item(8, 3).Value = "Fzd4"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.556150666666667
$ws.Cells.Item(8, 8).Value = 7.668452
$ws.Cells.Item(8, 9).Value = 0.4415485602524809
$ws.Cells.Item(8, 10).Value = 0.441548560252481
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.1305583333333333
$ws.Cells.Item(8, 14).Value = 0.391675
$ws.Cells.Item(8, 15).Value = 0.002064735720865253
$ws.Cells.Item(8, 16).Value = 0.002064735720865253
$ws.Cells.Item(8, 17).Value = 0.3337267707888889
$ws.Cells.Item(8, 18).Value = 3.0035409371
$ws.Cells.Item(8, 19).Value = 0.0009116810848499208
$ws.Cells.Item(8, 20).Value = 0.0009116810848499209

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Wnt9a"
$ws.Cells.Item(9, 3).Value = "Fzd4"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.556150666666667
$ws.Cells.Item(9, 8).Value = 7.668452
$ws.Cells.Item(9, 9).Value = 0.4415485602524809
$ws.Cells.Item(9, 10).Value = 0.441548560252481
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 23.275017
$ws.Cells.Item(9, 14).Value = 69.825051
$ws.Cells.Item(9, 15).Value = 0.3680864926557428
$ws.Cells.Item(9, 16).Value = 0.3680864926557428
$ws.Cells.Item(9, 17).Value = 59.49445022122801
$ws.Cells.Item(9, 18).Value = 535.4500519910521
$ws.Cells.Item(9, 19).Value = 0.1625280608805286
$ws.Cells.Item(9, 20).Value = 0.1625280608805286
